$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) cells remain text, matching the source format
# (values like "0.997" or "312.95" would otherwise be auto-converted to numbers).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '41.529.50'
$ws.Range('E2').Value = '  +0.71%  '
$ws.Range('D3').Value = '2.488.11'
$ws.Range('E3').Value = '  +1.17%  '
$ws.Range('D4').Value = '0.997'
$ws.Range('E4').Value = '  -0.26%  '
$ws.Range('D5').Value = '312.95'
$ws.Range('E5').Value = '  +0.28%  '
$ws.Range('D6').Value = '93.58'
$ws.Range('E6').Value = '  -0.72%  '
$ws.Range('D7').Value = '0.547'
$ws.Range('E7').Value = '  -0.84%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.26%  '
$ws.Range('D9').Value = '0.499'
$ws.Range('E9').Value = '  -0.68%  '
$ws.Range('D10').Value = '32.68'
$ws.Range('E10').Value = '  -2.29%  '
$ws.Range('D11').Value = '0.0783'
$ws.Range('E11').Value = '  +0.23%  '
$ws.Range('E12').Value = '  +1.98%  '
$ws.Range('D13').Value = '2.872.36'
$ws.Range('E13').Value = '  +1.12%  '
$ws.Range('D14').Value = '6.85'
$ws.Range('E14').Value = '  -1.54%  '
$ws.Range('D15').Value = '2.511.70'
$ws.Range('E15').Value = '  +2.32%  '
$ws.Range('D16').Value = '15.47'
$ws.Range('E16').Value = '  +6.00%  '
$ws.Range('D17').Value = '0.754'
$ws.Range('E17').Value = '  -4.27%  '
$ws.Range('D18').Value = '41.597.40'
$ws.Range('E18').Value = '  +0.99%  '
$ws.Range('D19').Value = '6.33'
$ws.Range('E19').Value = '  -0.37%  '
$ws.Range('E20').Value = '  +1.12%  '
$ws.Range('D21').Value = '70.69'
$ws.Range('E21').Value = '  +4.74%  '
$ws.Range('D22').Value = '11.21'
$ws.Range('E22').Value = '  -2.63%  '
$ws.Range('D23').Value = '235.87'
$ws.Range('E23').Value = '  -0.36%  '
$ws.Range('E24').Value = '  -2.58%  '
$ws.Range('B25').Value = 'ImmutableX'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D25').Value = '1.91'
$ws.Range('E25').Value = '  -1.17%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').Value = '24.93'
$ws.Range('E27').Value = '  +1.91%  '
$ws.Range('E28').Value = '  +0.35%  '
$ws.Range('D29').Value = '9.66'
$ws.Range('E29').Value = '  -0.43%  '
$ws.Range('D30').Value = '36.26'
$ws.Range('E30').Value = '  -0.26%  '
$ws.Range('D31').Value = '155.28'
$ws.Range('E31').Value = '  +1.53%  '
$ws.Range('D32').Value = '5.43'
$ws.Range('E32').Value = '  -2.87%  '
$ws.Range('B33').Value = 'Celestia'
$ws.Range('C33').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D33').Value = '18.24'
$ws.Range('E33').Value = '  +7.08%  '
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').Value = '2.57'
$ws.Range('E34').Value = '  -0.43%  '
$ws.Range('D35').Value = '0.0756'
$ws.Range('E35').Value = '  +0.50%  '
$ws.Range('E36').Value = '  -4.86%  '
$ws.Range('D37').Value = '2.96'
$ws.Range('E37').Value = '  -1.86%  '
$ws.Range('D38').Value = '1.85'
$ws.Range('E38').Value = '  -3.19%  '
$ws.Range('E39').Value = '  +1.56%  '
$ws.Range('E40').Value = '  -0.20%  '
$ws.Range('E41').Value = '  -3.39%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').Value = '20.28'
$ws.Range('E42').Value = '  -4.25%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  -0.17%  '
$ws.Range('D44').Value = '1.963.84'
$ws.Range('E44').Value = '  -0.19%  '
$ws.Range('E45').Value = '  -0.21%  '
$ws.Range('E46').Value = '  -3.50%  '
$ws.Range('E47').Value = '  +1.03%  '
$ws.Range('D48').Value = '2.729.38'
$ws.Range('E48').Value = '  +0.93%  '
$ws.Range('D49').Value = '96.55'
$ws.Range('E49').Value = '  -0.91%  '
$ws.Range('D50').Value = '67.40'
$ws.Range('E50').Value = '  -3.31%  '
$ws.Range('D51').Value = '73.36'
$ws.Range('E51').Value = '  -3.77%  '
